# JLCSMT_Sample_BOM1.xlsx -- "updated board layout to optimize for DFM"
#
# The BOM sheet carries a small "Board" placement/metadata range (the
# query-table refresh area, columns F:O) next to the actual BOM table
# (columns A:D). This edit:
#   1. shrinks the existing "Board" named range by one row (F1:O4 -> F1:O3)
#      and clears row 4 of that placeholder block (F4:J4)
#   2. adds a second "Board_1" named range (P1:Y7) for a newly added,
#      taller placeholder block (P1:T7)
#   3. re-splits the old row 4 BOM line (AW9523B / U2 / ... / C148077) into
#      four rows:
#        100n / C1,C2 / C0402 / C1525          (replaces old row 4 in place)
#        4k7 / R17,R18,R19 / R0402 / C25900    (new)
#        ATTINY2313A / U1 / QFN...21T260N / C185530   (new)
#        AW9523B / U2 / QFN...25T270N / C148077       (the old row 4, now row 7)
#   4. moves the active selection below the (now larger) table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: defined names -------------------------------------------------

$board = $wb.Names.Item("Board")
$board.RefersTo = "=Sheet1!`$F`$1:`$O`$3"

$ws.Names.Add("Board_1", "=Sheet1!`$P`$1:`$Y`$7")

# --- 3: BOM rows -----------------------------------------------------------

# Row 4 becomes the split-out ceramic cap line (was AW9523B/U2/.../C148077)
$ws.Range("A4").Value = "100n"
$ws.Range("B4").Value = "C1, C2"
$ws.Range("C4").Value = "C0402"
$ws.Range("D4").Value = "C1525"

# The F4:J4 placeholder cells belonged to the old (4-row) "Board" range;
# that range is now only 3 rows tall, so drop them entirely.
$ws.Range("F4:J4").Clear()

# New row 5
$ws.Range("A5").Value = "4k7"
$ws.Range("B5").Value = "R17, R18, R19"
$ws.Range("C5").Value = "R0402"
$ws.Range("D5").Value = "C25900"
$ws.Range("A5:C5").NumberFormat = "@"

# New row 6
$ws.Range("A6").Value = "ATTINY2313A"
$ws.Range("B6").Value = "U1"
$ws.Range("C6").Value = "QFN50P400X400X80-21T260N"
$ws.Range("D6").Value = "C185530"
$ws.Range("A6:C6").NumberFormat = "@"

# New row 7 -- the original row 4 contents, pushed down
$ws.Range("A7").Value = "AW9523B"
$ws.Range("B7").Value = "U2"
$ws.Range("C7").Value = "QFN50P400X400X80-25T270N"
$ws.Range("D7").Value = "C148077"
$ws.Range("A7:C7").NumberFormat = "@"

# --- 4: new P:T placeholder block (Board_1), rows 1-7 ----------------------

$ws.Range("P1:T1").NumberFormat = "@"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("R2:T2").NumberFormat = "@"
$ws.Range("P3:T3").NumberFormat = "@"
$ws.Range("P4:T4").NumberFormat = "@"
$ws.Range("P5:T5").NumberFormat = "@"
$ws.Range("P6:T6").NumberFormat = "@"
$ws.Range("P7:T7").NumberFormat = "@"

$ws.Columns("P").ColumnWidth = $ws.Columns("F").ColumnWidth
$ws.Columns("Q:S").ColumnWidth = $ws.Columns("G").ColumnWidth
$ws.Columns("T").ColumnWidth = $ws.Columns("J").ColumnWidth
$ws.Columns("U").ColumnWidth = $ws.Columns("K").ColumnWidth
$ws.Columns("V").ColumnWidth = $ws.Columns("L").ColumnWidth
$ws.Columns("W:X").ColumnWidth = $ws.Columns("M").ColumnWidth

# --- 5: selection below the table ------------------------------------------

$ws.Range("A8:XFD12").Select()
